$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 278, pushing the existing rows 278:291 down to 279:292.
$ws.Rows.Item(278).Insert()

# Populate the newly inserted row 278 with the new weekly record.
$ws.Cells.Item(278, 1).Value2 = 6
$ws.Cells.Item(278, 2).Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(278, 3).Value2 = "Metropolitana"
$ws.Cells.Item(278, 4).Value2 = 44585
$ws.Cells.Item(278, 5).Value2 = 13
$ws.Cells.Item(278, 6).Value2 = "Fruta"
$ws.Cells.Item(278, 7).Value2 = 100101
$ws.Cells.Item(278, 8).Value2 = "Berries"
$ws.Cells.Item(278, 9).Value2 = 100101001
$ws.Cells.Item(278, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(278, 11).Value2 = "Sin especificar"
$ws.Cells.Item(278, 12).Value2 = "Primera"
$ws.Cells.Item(278, 13).Value2 = 2000
$ws.Cells.Item(278, 14).Value2 = 4000
$ws.Cells.Item(278, 15).Value2 = 4000
$ws.Cells.Item(278, 16).Value2 = 4000
$ws.Cells.Item(278, 17).Value2 = "`$/bandeja 2 kilos"
$ws.Cells.Item(278, 18).Value2 = "Provincia de Linares"
$ws.Cells.Item(278, 19).Value2 = 2000
$ws.Cells.Item(278, 20).Value2 = 2
